# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps for the second
# data row (row 3) on both the "zh-cn" and "de-de" report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 14:19:46"
$wsZhCn.Range("G3").Value = "2016-01-08 14:20:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 14:20:00"
$wsDeDe.Range("G3").Value = "2016-01-08 14:20:56"
